# Fruta / hortaliza, semanal
# Insert a new weekly price-report row above the current row 233 (pushes the
# existing rows 233:247 down to 234:248) and populate it with this week's
# data for "Uva" / "Red Globe" at Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 233:247 down to 234:248, leaving a blank row 233 (formatted
# like the row above it, same as Excel's native Insert behaviour).
$ws.Rows.Item(233).Insert()

# Fill in the new row with the latest week's figures.
$ws.Cells.Item(233, 1).Value = 4
$ws.Cells.Item(233, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(233, 3).Value = "Los Lagos"
$ws.Cells.Item(233, 4).Value = 44746
$ws.Cells.Item(233, 5).Value = 10
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100109
$ws.Cells.Item(233, 8).Value = "Uva"
$ws.Cells.Item(233, 9).Value = 100109001
$ws.Cells.Item(233, 10).Value = "Uva"
$ws.Cells.Item(233, 11).Value = "Red Globe"
$ws.Cells.Item(233, 12).Value = "Primera"
$ws.Cells.Item(233, 13).Value = 200
$ws.Cells.Item(233, 14).Value = 8000
$ws.Cells.Item(233, 15).Value = 9000
$ws.Cells.Item(233, 16).Value = 8500
$ws.Cells.Item(233, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(233, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(233, 19).Value = 425
$ws.Cells.Item(233, 20).Value = 20
